$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "data as of" timestamp
$ws.Range("A1").Value = "Datos actualizados a 6 de Abril de 2020 a las 15:22"

# Update country case statistics (refreshed counts; a few countries swap rank
# order, which is reflected by swapping the country-name cells in column A)
$ws.Range("B4").Value = 336906
$ws.Range("C4").Value = 233
$ws.Range("E4").Value = 309305
$ws.Range("G4").Value = 8
$ws.Range("H4").Value = 9624
$ws.Range("B7").Value = 100232
$ws.Range("C7").Value = 109
$ws.Range("E7").Value = 69941
$ws.Range("G7").Value = 7
$ws.Range("H7").Value = 1591
$ws.Range("B17").Value = 12206
$ws.Range("C17").Value = 155
$ws.Range("E17").Value = 8523
$ws.Range("E21").Value = 7971
$ws.Range("G21").Value = 6
$ws.Range("H21").Value = 55
$ws.Range("E32").Value = 3489
$ws.Range("G32").Value = 11
$ws.Range("H32").Value = 162
$ws.Range("A39").Value = "Arabia Saudita"
$ws.Range("B39").Value = 2523
$ws.Range("C39").Value = 121
$ws.Range("D39").Value = 551
$ws.Range("E39").Value = 1934
$ws.Range("F39").Value = 41
$ws.Range("G39").Value = 4
$ws.Range("H39").Value = 38
$ws.Range("A40").Value = "Indonesia"
$ws.Range("B40").Value = 2491
$ws.Range("C40").Value = 218
$ws.Range("D40").Value = 192
$ws.Range("E40").Value = 2090
$ws.Range("F40").Value = 0
$ws.Range("G40").Value = 11
$ws.Range("H40").Value = 209
$ws.Range("E43").Value = 1849
$ws.Range("F43").Value = 81
$ws.Range("G43").Value = 1
$ws.Range("H43").Value = 27
$ws.Range("A45").Value = "Emiratos Arabes Unidos"
$ws.Range("B45").Value = 2076
$ws.Range("C45").Value = 277
$ws.Range("D45").Value = 167
$ws.Range("E45").Value = 1898
$ws.Range("F45").Value = 1
$ws.Range("G45").Value = 1
$ws.Range("H45").Value = 11
$ws.Range("A46").Value = "Panama"
$ws.Range("B46").Value = 1988
$ws.Range("C46").Value = 187
$ws.Range("D46").Value = 13
$ws.Range("E46").Value = 1921
$ws.Range("F46").Value = 78
$ws.Range("G46").Value = 8
$ws.Range("H46").Value = 54
$ws.Range("A47").Value = "Serbia"
$ws.Range("B47").Value = 1908
$ws.Range("D47").Value = 54
$ws.Range("E47").Value = 1803
$ws.Range("F47").Value = 98
$ws.Range("H47").Value = 51
$ws.Range("E49").Value = 1581
$ws.Range("G49").Value = 3
$ws.Range("H49").Value = 76
$ws.Range("A52").Value = "Islandia"
$ws.Range("B52").Value = 1562
$ws.Range("C52").Value = 76
$ws.Range("D52").Value = 460
$ws.Range("E52").Value = 1098
$ws.Range("F52").Value = 12
$ws.Range("G52").Value = 0
$ws.Range("H52").Value = 4
$ws.Range("A53").Value = "Argentina"
$ws.Range("B53").Value = 1554
$ws.Range("D53").Value = 280
$ws.Range("E53").Value = 1226
$ws.Range("F53").Value = 86
$ws.Range("G53").Value = 2
$ws.Range("H53").Value = 48
$ws.Range("B73").Value = 674
$ws.Range("C73").Value = 20
$ws.Range("E73").Value = 598
$ws.Range("D113").Value = 38
$ws.Range("E113").Value = 135
$ws.Range("A200").Value = "Malaui"
$ws.Range("C200").Value = 1
$ws.Range("E200").Value = 5
$ws.Range("G200").Value = 0
$ws.Range("H200").Value = 0
$ws.Range("A201").Value = "Belice"
$ws.Range("F201").Value = 1
$ws.Range("G201").Value = 1
$ws.Range("A202").Value = "Islas Turcas y Caicos"
$ws.Range("D202").Value = 0
$ws.Range("E202").Value = 4
$ws.Range("H202").Value = 1
$ws.Range("A203").Value = "Butan"
$ws.Range("B203").Value = 5
$ws.Range("D203").Value = 2
$ws.Range("E203").Value = 3
$ws.Range("A204").Value = "Sahara Occidental"
$ws.Range("A206").Value = "Burundi"
$ws.Range("A207").Value = "Islas Virgenes Britanicas"
$ws.Range("A208").Value = "Anguila"
$ws.Range("A209").Value = "Bonaire, San Eustaquio y Saba"
$ws.Range("A210").Value = "Papua Nueva Guinea"
$ws.Range("C210").Value = 1
$ws.Range("A211").Value = "Islas Malvinas"
$ws.Range("C211").Value = 0
